$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","E","F","G","H","I","J","K","L","M")

$row8 = @("9 ماهه منتهی به 1399/09","12 ماهه منتهی به 1399/12","3 ماهه منتهی به 1400/03","6 ماهه منتهی به 1400/06","9 ماهه منتهی به 1400/09","12 ماهه منتهی به 1400/12","3 ماهه منتهی به 1401/03","6 ماهه منتهی به 1401/06","9 ماهه منتهی به 1401/09","12 ماهه منتهی به 1401/12")
$row9 = @("1400-10-29 (2)","1401-02-06 (9)","1401-04-29 (3)","1401-08-25 (4)","1401-10-28 (2)","1402-02-10 (8)","1401-04-29","1401-08-25 (2)","1401-10-28","1402-02-10 (2)")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $row8[$i]
    $ws.Range($cols[$i] + "9").Value = $row9[$i]
}

$row11 = @(2389378,3409060,1721437,3777997,6073241,8310129,3138299,5925440,8601081,11031555)
$row12 = @(-1753818,-2484486,-1266016,-2814649,-4645370,-6576671,-2484291,-4654459,-6805429,-8781736)
$row13 = @(635560,924574,455421,963348,1427871,1733458,654008,1270981,1795652,2249819)
$row14 = @(-74122,-117229,-32002,-71302,-103850,-163847,-45953,-88820,-138691,-215884)
$row15 = @(0,0,0,0,0,0,0,0,0,0)
$row16 = @(7365,30553,25345,57916,88792,115502,19000,35304,82075,33125)
$row17 = @(568803,837898,448764,949962,1412813,1685113,627055,1217465,1739036,2067060)
$row18 = @(-8823,-9111,-3048,-8228,-13461,-14203,-4265,-4376,-5319,-13393)
$row19 = @(10946,-44015,210,17187,21437,17174,1928,9177,19269,-12868)
$row20 = @(570926,784772,445926,958921,1420789,1688084,624718,1222266,1752986,2040799)
$row21 = @(-112198,-113671,-100333,-168000,-282850,-282485,-121821,-188155,-298008,-209718)
$row22 = @(458728,671101,345593,790921,1137939,1405599,502897,1034111,1454978,1831081)
$row23 = @(0,0,0,0,0,0,0,0,0,0)
$row24 = @(458728,671101,345593,790921,1137939,1405599,502897,1034111,1454978,1831081)
$row25 = @(967,1414,728,1667,2398,2962,1060,1089,1532,1927)
$row26 = @(474522,474522,474522,474522,474522,474522,474522,950000,950000,950000)
$row27 = @(483,706,364,833,1198,1480,529,1089,1532,1927)

$rows = @(11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27)
$rowVals = @($row11,$row12,$row13,$row14,$row15,$row16,$row17,$row18,$row19,$row20,$row21,$row22,$row23,$row24,$row25,$row26,$row27)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $vals = $rowVals[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $rows[$r]).Value = $vals[$i]
    }
}

$ws.Columns("D").ColumnWidth = 27.166666666666668
$ws.Columns("E").ColumnWidth = 28.166666666666668
$ws.Columns("F").ColumnWidth = 27.166666666666668
$ws.Columns("G").ColumnWidth = 27.166666666666668
$ws.Columns("H").ColumnWidth = 27.166666666666668
$ws.Columns("I").ColumnWidth = 28.166666666666668
$ws.Columns("J").ColumnWidth = 27.166666666666668
$ws.Columns("K").ColumnWidth = 27.166666666666668
$ws.Columns("L").ColumnWidth = 27.166666666666668
$ws.Columns("M").ColumnWidth = 28.166666666666668